$d = $word.ActiveDocument

# 1. Heading: "Danh " + "sách" + " " + "thành" + " " + "viên" -> single run "Danh sách thành viên"
$d.Content.Find.Execute("Danh sách thành viên", $true, $false, $false, $false, $false, $true, 1, $false, "Danh sách thành viên", 2)

# 2. Table header: "Họ" + " " + "Tên" -> "Họ Tên"
$d.Content.Find.Execute("Họ Tên", $true, $false, $false, $false, $false, $true, 1, $false, "Họ Tên", 2)

# 3. "Nguyễn" + " Đức " + "Hòa" -> "Nguyễn Đức Hòa"
$d.Content.Find.Execute("Nguyễn Đức Hòa", $true, $false, $false, $false, $false, $true, 1, $false, "Nguyễn Đức Hòa", 2)

# 4. "Nguyễn" + " Võ " + "Nguyên" + " " + "Phát" -> "Nguyễn Võ Nguyên Phát"
$d.Content.Find.Execute("Nguyễn Võ Nguyên Phát", $true, $false, $false, $false, $false, $true, 1, $false, "Nguyễn Võ Nguyên Phát", 2)

# 5. MSSV fix: 21242042 -> 21424042
$d.Content.Find.Execute("21242042", $true, $false, $false, $false, $false, $true, 1, $false, "21424042", 2)

# 6. "Trần" + " Văn " + "Trụ" + " (" + "Người" + " quay clip " + "nên" + " " + "không" + " " + "thấy" + " " + "tên" + " " + "trên" + " video)"
#    -> "Trần Văn Trụ (Người quay clip nên không thấy tên trên video)"
$d.Content.Find.Execute("Trần Văn Trụ (Người quay clip nên không thấy tên trên video)", $true, $false, $false, $false, $false, $true, 1, $false, "Trần Văn Trụ (Người quay clip nên không thấy tên trên video)", 2)

# 7. "Phạm Minh " + "Toàn" -> "Phạm Minh Toàn"
$d.Content.Find.Execute("Phạm Minh Toàn", $true, $false, $false, $false, $false, $true, 1, $false, "Phạm Minh Toàn", 2)

# 8. "Nguyễn" + " Đức Thịnh" -> "Nguyễn Đức Thịnh"
$d.Content.Find.Execute("Nguyễn Đức Thịnh", $true, $false, $false, $false, $false, $true, 1, $false, "Nguyễn Đức Thịnh", 2)

# 9. "Link" + " google drive:" -> "Link google drive:"
$d.Content.Find.Execute("Link google drive:", $true, $false, $false, $false, $false, $true, 1, $false, "Link google drive:", 2)
